$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 ("Rejection reasons") has its populated values shifted one column to
# the right (D->E, H->I, K->L, R->S, T->U, W->X, X->Y, AA->AB, AG->AH),
# and a new value "23" inserted at F21.
# Process from right to left so sources aren't overwritten before being read.

$ws.Range("AH21").Value2 = $ws.Range("AG21").Value2
$ws.Range("AG21").Value2 = $null

$ws.Range("AB21").Value2 = $ws.Range("AA21").Value2
$ws.Range("AA21").Value2 = $null

$ws.Range("Y21").Value2 = $ws.Range("X21").Value2
$ws.Range("X21").Value2 = $ws.Range("W21").Value2
$ws.Range("W21").Value2 = $null

$ws.Range("U21").Value2 = $ws.Range("T21").Value2
$ws.Range("T21").Value2 = $null

$ws.Range("S21").Value2 = $ws.Range("R21").Value2
$ws.Range("R21").Value2 = $null

$ws.Range("L21").Value2 = $ws.Range("K21").Value2
$ws.Range("K21").Value2 = $null

$ws.Range("I21").Value2 = $ws.Range("H21").Value2
$ws.Range("H21").Value2 = $null

$ws.Range("E21").Value2 = $ws.Range("D21").Value2
$ws.Range("D21").Value2 = $null

$ws.Range("F21").Value2 = "23"
